# Update date line and all division exercise answers per commit diff
$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-07-05 Saturday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-07-06 Sunday", 2) | Out-Null
$d.Content.Find.Execute("63÷7=9, 0", $true, $false, $false, $false, $false, $true, 1, $false, "82÷6=13, 4", 2) | Out-Null
$d.Content.Find.Execute("79÷7=11, 2", $true, $false, $false, $false, $false, $true, 1, $false, "53÷4=13, 1", 2) | Out-Null
$d.Content.Find.Execute("39÷5=7, 4", $true, $false, $false, $false, $false, $true, 1, $false, "60÷7=8, 4", 2) | Out-Null
$d.Content.Find.Execute("51÷4=12, 3", $true, $false, $false, $false, $false, $true, 1, $false, "76÷2=38, 0", 2) | Out-Null
$d.Content.Find.Execute("47÷6=7, 5", $true, $false, $false, $false, $false, $true, 1, $false, "61÷2=30, 1", 2) | Out-Null
$d.Content.Find.Execute("93÷8=11, 5", $true, $false, $false, $false, $false, $true, 1, $false, "97÷3=32, 1", 2) | Out-Null
$d.Content.Find.Execute("48÷7=6, 6", $true, $false, $false, $false, $false, $true, 1, $false, "17÷7=2, 3", 2) | Out-Null
$d.Content.Find.Execute("27÷8=3, 3", $true, $false, $false, $false, $false, $true, 1, $false, "89÷4=22, 1", 2) | Out-Null
$d.Content.Find.Execute("18÷3=6, 0", $true, $false, $false, $false, $false, $true, 1, $false, "96÷6=16, 0", 2) | Out-Null
$d.Content.Find.Execute("47÷3=15, 2", $true, $false, $false, $false, $false, $true, 1, $false, "51÷2=25, 1", 2) | Out-Null
$d.Content.Find.Execute("16÷6=2, 4", $true, $false, $false, $false, $false, $true, 1, $false, "40÷4=10, 0", 2) | Out-Null
$d.Content.Find.Execute("65÷8=8, 1", $true, $false, $false, $false, $false, $true, 1, $false, "90÷7=12, 6", 2) | Out-Null
$d.Content.Find.Execute("31÷2=15, 1", $true, $false, $false, $false, $false, $true, 1, $false, "21÷2=10, 1", 2) | Out-Null
$d.Content.Find.Execute("69÷3=23, 0", $true, $false, $false, $false, $false, $true, 1, $false, "67÷7=9, 4", 2) | Out-Null
$d.Content.Find.Execute("49÷4=12, 1", $true, $false, $false, $false, $false, $true, 1, $false, "32÷2=16, 0", 2) | Out-Null
$d.Content.Find.Execute("80÷3=26, 2", $true, $false, $false, $false, $false, $true, 1, $false, "54÷2=27, 0", 2) | Out-Null
$d.Content.Find.Execute("14÷2=7, 0", $true, $false, $false, $false, $false, $true, 1, $false, "47÷6=7, 5", 2) | Out-Null
$d.Content.Find.Execute("99÷4=24, 3", $true, $false, $false, $false, $false, $true, 1, $false, "70÷2=35, 0", 2) | Out-Null
$d.Content.Find.Execute("26÷2=13, 0", $true, $false, $false, $false, $false, $true, 1, $false, "86÷7=12, 2", 2) | Out-Null
$d.Content.Find.Execute("53÷3=17, 2", $true, $false, $false, $false, $false, $true, 1, $false, "25÷4=6, 1", 2) | Out-Null
$d.Content.Find.Execute("26÷5=5, 1", $true, $false, $false, $false, $false, $true, 1, $false, "62÷9=6, 8", 2) | Out-Null
$d.Content.Find.Execute("11÷8=1, 3", $true, $false, $false, $false, $false, $true, 1, $false, "30÷9=3, 3", 2) | Out-Null
$d.Content.Find.Execute("68÷4=17, 0", $true, $false, $false, $false, $false, $true, 1, $false, "18÷2=9, 0", 2) | Out-Null
$d.Content.Find.Execute("64÷8=8, 0", $true, $false, $false, $false, $false, $true, 1, $false, "85÷4=21, 1", 2) | Out-Null
$d.Content.Find.Execute("21÷8=2, 5", $true, $false, $false, $false, $false, $true, 1, $false, "67÷3=22, 1", 2) | Out-Null

Write-Output "Replacements complete"
